$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from an existing header cell (AC1) onto the new
# header cells so they keep the same bold/centered/bordered look.
$headerStyleSrc = $ws.Range("AC1")
$headerStyleSrc.Copy()
$newHeaders = $ws.Range("AD1:AF1")
$newHeaders.PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record (Wins/Losses/Ties) for every player row.
$lastRow = 62
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 67
    $ws.Cells.Item($r, 31).Value = 95
    $ws.Cells.Item($r, 32).Value = 0
}

$excel.CutCopyMode = 0
